$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data in column A (A1:A6)
$ws.Range("A1:A6").ClearContents()

# Write header into B1
$ws.Range("B1").Value = "house_size"

# Write the new data values into B2:B7
$values = @(2433, 2503, 3300, 3203, 3300, 4003)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Update the selection to match the target state
$ws.Range("E10").Select()
